# CE09OSPM added eng bar codes and bogey for D4 recovery date
#
# 1. Asset_Cal_Info!E31 (the STCENG000 / OSPM-00001-STC row) gets its
#    Sensor OOIBARCODE populated with the newly labeled engineering bar
#    code "OL000372" (was blank).
# 2. The placeholder "bogey" row for the RTE000000 asset (row 32, with a
#    bare serial number "950" and no bar code) is removed entirely.
# 3. The legacy AutoFilter-derived defined names that pointed at the old
#    398-row-plus-header range ($A$1:$H$399) are updated to reflect the
#    shrunk table ($A$1:$H$398) now that a row has been deleted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

# Add the engineering bar code for the STCENG000 sensor row.
$ws.Range("E31").Value = "OL000372"
$ws.Range("E31").ClearFormats()

# Remove the bogey RTE000000 placeholder row entirely (shifts rows below up).
$ws.Rows.Item(32).Delete()

# Keep the stale AutoFilter-derived named ranges in sync with the new
# (one-row-shorter) table extent.
$namesToFix = @(
    "_FilterDatabase_0_0_0_0_0_0",
    "_FilterDatabase_0_0_0_0_0_0_0_0",
    "_FilterDatabase_0_0_0_0_1",
    "_FilterDatabase_0_0_1",
    "_FilterDatabase_2"
)
foreach ($name in $namesToFix) {
    $n = $wb.Names.Item($name)
    $n.RefersTo = "=Asset_Cal_Info!`$A`$1:`$H`$398"
}
